$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to reflect the latest scrape.
# Numeric-looking price strings are prefixed with a leading apostrophe so
# Excel stores them as text (preserving formatting like trailing zeros and
# thousand-separator dots), matching how the source data is represented.
$ws.Range("D2").Value = "56.888.59"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "2.505.40"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D5").Value = "'496.70"
$ws.Range("E5").Value = "  +3.22%  "
$ws.Range("D6").Value = "'153.74"
$ws.Range("E6").Value = "  +9.79%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").Value = "2.514.51"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "'5.78"
$ws.Range("E10").Value = "  +6.00%  "
$ws.Range("D11").Value = "'0.0994"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "'0.337"
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "2.943.89"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "56.941.41"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").Value = "'21.47"
$ws.Range("E16").Value = "  +4.20%  "
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "2.519.92"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "'4.57"
$ws.Range("E19").Value = "  +4.41%  "
$ws.Range("E20").Value = "  +3.45%  "
$ws.Range("D21").Value = "'324.66"
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("D22").Value = "'0.996"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'5.94"
$ws.Range("E23").Value = "  +4.42%  "
$ws.Range("D24").Value = "'59.12"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("D25").Value = "'0.413"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").Value = "'0.164"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "2.608.87"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").Value = "'7.73"
$ws.Range("E29").Value = "  +4.38%  "
$ws.Range("D30").Value = "0.0₃0820"
$ws.Range("E30").Value = "  +4.91%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'151.38"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").Value = "'18.47"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("E36").Value = "  +5.67%  "
$ws.Range("E37").Value = "  +3.31%  "
$ws.Range("D38").Value = "'0.886"
$ws.Range("E38").Value = "  +4.86%  "
$ws.Range("E39").Value = "  +6.00%  "
$ws.Range("D40").Value = "'34.23"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.53"
$ws.Range("E41").Value = "  +3.97%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0567"
$ws.Range("E42").Value = "  +3.22%  "
$ws.Range("D43").Value = "'0.615"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").Value = "'0.994"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "'4.97"
$ws.Range("E45").Value = "  +10.38%  "
$ws.Range("D46").Value = "'271.26"
$ws.Range("E46").Value = "  +8.84%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0928"
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0232"
$ws.Range("E48").Value = "  +3.69%  "
$ws.Range("D49").Value = "'10.20"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "'18.05"
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("D51").Value = "1.914.72"
$ws.Range("E51").Value = "  -2.87%  "
